$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "player" / "Player" translation pair right before
# the "When Activating/Deactivating/Inverting" event rows (old row 35),
# shifting those rows down by one (35->36, 36->37, 37->38).
$ws.Rows(34).Insert()

$ws.Range("A34").Value = "player"
$ws.Range("B34").Value = "Player"

# The shifted event rows previously carried the small (10pt) font style;
# reset them back to the Normal/default look.
$ws.Range("A36:B38").Style = "Normal"

# Restore view state (zoom / scroll / selection) similar to the saved file.
$excel.ActiveWindow.Zoom = 115
$ws.Range("C34").Select()
